$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary block (rows 14-17): label in column A, stat in column B ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Copy B14's formatting (bold, size 12, vertical-centered) onto the rest of the block
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the new labeled rows
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Average of |S*|/n in row 12 (bold, default size)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Selection left active over the new block, as in the saved file
$ws.Range("A14:B17").Select() | Out-Null

# Page setup (paper size A4 / portrait) matching the resaved workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
